# Update betting-odds values on Sheet1 to reflect the refreshed
# FlashScore odds snapshot (see commit "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Corinthians x Athletico-PR)
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.97
$ws.Range("R2").Value = 1.93

# Row 3 (Flamengo RJ x Fluminense)
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 6.5

# Row 5 (Shrewsbury x Exeter)
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("X5").Value = 13
$ws.Range("AN5").Value = 4.75
$ws.Range("AO5").Value = 15

# Row 7 (Atl. Morelia x Zacatecas Mineros)
$ws.Range("G7").Value = 2.87
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 2.42
$ws.Range("J7").Value = 3.35
$ws.Range("K7").Value = 2.07
$ws.Range("L7").Value = 2.95
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 7.7
$ws.Range("O7").Value = 1.35
$ws.Range("P7").Value = 2.7
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.62
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.83
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 10.5
$ws.Range("Z7").Value = 35
$ws.Range("AA7").Value = 26
$ws.Range("AB7").Value = 35
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 5.9
$ws.Range("AE7").Value = 14
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 600
$ws.Range("AH7").Value = 7.2
$ws.Range("AI7").Value = 11.5
$ws.Range("AM7").Value = 32
$ws.Range("AN7").Value = 4.8
$ws.Range("AO7").Value = 15
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 70
$ws.Range("AR7").Value = 90
$ws.Range("AT7").Value = 2.55
$ws.Range("AU7").Value = 6.6
$ws.Range("AV7").Value = 55
$ws.Range("AW7").Value = 4.3
$ws.Range("AY7").Value = 19.5
$ws.Range("AZ7").Value = 50
$ws.Range("BA7").Value = 80

# Row 9 (Cusco x Los Chankas)
$ws.Range("O9").Value = 1.14
$ws.Range("P9").Value = 5.5

# Row 12 (Fenix x CA Cerro)
$ws.Range("G12").Value = 2.1
$ws.Range("I12").Value = 3.75
$ws.Range("K12").Value = 1.91
$ws.Range("Q12").Value = 2.6
$ws.Range("R12").Value = 1.48
$ws.Range("S12").Value = 1.57
$ws.Range("T12").Value = 2.25
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 1.62
$ws.Range("W12").Value = 5.5
$ws.Range("Y12").Value = 10
$ws.Range("AE12").Value = 21
$ws.Range("AF12").Value = 81
$ws.Range("AH12").Value = 8
$ws.Range("AP12").Value = 29
$ws.Range("AT12").Value = 2.25
$ws.Range("AU12").Value = 9.5
$ws.Range("AV12").Value = 81
$ws.Range("AY12").Value = 41
